$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Last Updated timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(2, 1).Value = "30 Oct 2025, 10:45 AM"

# --- "1 Month Performance" sheet: updated stock rankings/percentages ---
$perf = $wb.Worksheets.Item("1 Month Performance")
$perf.Cells.Item(4, 3).Value = 82.5116
$perf.Cells.Item(5, 3).Value = 69.7568
$perf.Cells.Item(7, 3).Value = 57.4261
$perf.Cells.Item(8, 3).Value = 57.289
$perf.Cells.Item(9, 3).Value = 52.0057
$perf.Cells.Item(12, 3).Value = 40.4496
$perf.Cells.Item(13, 3).Value = 37.7864
$perf.Cells.Item(15, 3).Value = 37.2631
$perf.Cells.Item(16, 3).Value = 37.2587
$perf.Cells.Item(17, 3).Value = 36.6144
$perf.Cells.Item(18, 3).Value = 36.3616
$perf.Cells.Item(20, 2).Value = "MAANALU"
$perf.Cells.Item(20, 3).Value = 33.6283
$perf.Cells.Item(21, 2).Value = "RAMAPHO"
$perf.Cells.Item(21, 3).Value = 33.5979
$perf.Cells.Item(22, 2).Value = "MEGASOFT"
$perf.Cells.Item(22, 3).Value = 33.4275
$perf.Cells.Item(23, 2).Value = "SOUTHBANK"
$perf.Cells.Item(23, 3).Value = 33.3103
$perf.Cells.Item(25, 2).Value = "ORIENTTECH"
$perf.Cells.Item(25, 3).Value = 29.9312
$perf.Cells.Item(26, 2).Value = "MRPL"
$perf.Cells.Item(26, 3).Value = 29.749
$perf.Cells.Item(27, 2).Value = "INDORAMA"
$perf.Cells.Item(27, 3).Value = 28.6141
$perf.Cells.Item(28, 2).Value = "ARFIN"
$perf.Cells.Item(28, 3).Value = 28.5781
$perf.Cells.Item(30, 3).Value = 27.0325
$perf.Cells.Item(32, 3).Value = 26.5004
$perf.Cells.Item(33, 2).Value = "RAMCOSYS"
$perf.Cells.Item(33, 3).Value = 25.8359
$perf.Cells.Item(34, 2).Value = "HATSUN"
$perf.Cells.Item(34, 3).Value = 25.6407
$perf.Cells.Item(35, 2).Value = "SKYGOLD"
$perf.Cells.Item(35, 3).Value = 25.1596
$perf.Cells.Item(36, 2).Value = "SAGILITY"
$perf.Cells.Item(36, 3).Value = 25.0615
$perf.Cells.Item(37, 3).Value = 24.8787
$perf.Cells.Item(39, 3).Value = 24.499
$perf.Cells.Item(41, 3).Value = 24.0157
$perf.Cells.Item(44, 2).Value = "INDRAMEDCO"
$perf.Cells.Item(44, 3).Value = 23.7469
$perf.Cells.Item(45, 2).Value = "KERNEX"
$perf.Cells.Item(45, 3).Value = 23.5234
$perf.Cells.Item(46, 3).Value = 23.4583
$perf.Cells.Item(47, 3).Value = 23.4149
$perf.Cells.Item(48, 2).Value = "LORDSCHLO"
$perf.Cells.Item(48, 3).Value = 22.7904
$perf.Cells.Item(49, 2).Value = "DCBBANK"
$perf.Cells.Item(49, 3).Value = 22.6612
$perf.Cells.Item(50, 3).Value = 21.9276
$perf.Cells.Item(52, 3).Value = 21.0891
$perf.Cells.Item(53, 3).Value = 21.013
$perf.Cells.Item(54, 2).Value = "GUJTHEM"
$perf.Cells.Item(54, 3).Value = 20.8787
$perf.Cells.Item(55, 2).Value = "MOLDTECH"
$perf.Cells.Item(55, 3).Value = 20.8659
$perf.Cells.Item(56, 2).Value = "RBLBANK"
$perf.Cells.Item(56, 3).Value = 20.8567
$perf.Cells.Item(57, 2).Value = "MARINE"
$perf.Cells.Item(57, 3).Value = 20.8447
$perf.Cells.Item(58, 2).Value = "SKMEGGPROD"
$perf.Cells.Item(58, 3).Value = 20.7435
$perf.Cells.Item(59, 2).Value = "BHARATWIRE"
$perf.Cells.Item(59, 3).Value = 20.7357
$perf.Cells.Item(60, 2).Value = "SCI"
$perf.Cells.Item(60, 3).Value = 20.7102
$perf.Cells.Item(61, 2).Value = "FEDERALBNK"
$perf.Cells.Item(61, 3).Value = 20.2179
$perf.Cells.Item(62, 2).Value = "UNIPARTS"
$perf.Cells.Item(62, 3).Value = 20.1592
$perf.Cells.Item(63, 3).Value = 20.0858
$perf.Cells.Item(64, 3).Value = 19.8307
$perf.Cells.Item(65, 3).Value = 19.6664
$perf.Cells.Item(66, 3).Value = 19.6222
$perf.Cells.Item(67, 2).Value = "BANKINDIA"
$perf.Cells.Item(67, 3).Value = 19.578
$perf.Cells.Item(68, 2).Value = "REPRO"
$perf.Cells.Item(68, 3).Value = 19.5397
$perf.Cells.Item(70, 3).Value = 19.3382
$perf.Cells.Item(71, 3).Value = 19.1294
$perf.Cells.Item(72, 3).Value = 19.0845
$perf.Cells.Item(73, 3).Value = 18.57
$perf.Cells.Item(74, 2).Value = "THOMASCOTT"
$perf.Cells.Item(74, 3).Value = 18.3343
$perf.Cells.Item(76, 2).Value = "ABDL"
$perf.Cells.Item(76, 3).Value = 18.1373
